$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D ("canonical SMILES")
$ws.Range("D2").Value2 = "canonical SMILES"

# Duplicate the "canonical isomeric SMILES" values (column C) into the new
# "canonical SMILES" column (column D) for every data row.
for ($r = 3; $r -le 10; $r++) {
    $ws.Cells.Item($r, 4).Value2 = $ws.Cells.Item($r, 3).Value2
}

# Size the new column like the others (renders as width="37" in the xlsx).
$ws.Columns.Item(4).ColumnWidth = 36.15
